$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "SMGAN" column (B). Excel shifts the
# existing column B (SMGAN) formatting/content into column C automatically,
# and gives the freshly inserted column B a copy of column A's cell styles.
$ws.Columns("B").Insert()

# --- Header row ---------------------------------------------------------
# B1 keeps holding the (renamed) "SMGAN" shared string -> becomes "SMGAN BASELINE"
$ws.Range("B1").Value = "SMGAN BASELINE"
# New column C1 gets the brand new header
$ws.Range("C1").Value = "SMGAN Tuned"
# D1 (old C1, "MolGAN") is untouched by the insert/shift.

# --- New "SMGAN BASELINE" data column (B) -------------------------------
$ws.Range("B2").Value = 17.5
$ws.Range("B3").Value = 0.02
$ws.Range("B4").Value = 100
$ws.Range("B5").Value = 0.58
$ws.Range("B6").Value = 0.03
$ws.Range("B7").Value = 0.42

# --- Column widths -------------------------------------------------------
# Match column B's width to column A's, and nudge column C close to its
# auto-fit width for the new "SMGAN Tuned" header text.
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth
$ws.Columns("C").ColumnWidth = 13.25

$ws.Range("A1:D7").Select() | Out-Null
